$wb = $excel.ActiveWorkbook

# --- Rename the current "Graduate Students" sheet to "Full-time Graduate Students" ---
$wsFullTime = $wb.Worksheets.Item("Graduate Students")
$wsFullTime.Name = "Full-time Graduate Students"

# --- Insert a new worksheet right after it for the part-time data ---
$wsPartTime = $wb.Worksheets.Add($null, $wsFullTime)
$wsPartTime.Name = "Part-time Graduate Students"
$wsPartTime.StandardHeight = 16

$ptLabels = @(
    "All races"
    ,"All part-time students"
    ,"    Male"
    ,"    Female"
    ,"  U.S. citizens and permanent residents"
    ,"    Hispanic or Latino"
    ,"    Not Hispanic or Latino"
    ,"      American Indian or Alaska Native"
    ,"      Asian"
    ,"      Black or African American"
    ,"      Native Hawaiian or Other Pacific Islander"
    ,"      White"
    ,"      More than one race"
    ,"    Unknown ethnicity and race"
    ,"  Temporary visa holders"
    ,"  Science and engineering"
    ,"    Science"
    ,"      Agricultural and veterinary sciences"
    ,"      Biological and biomedical sciences"
    ,"      Communication"
    ,"      Computer and information sciences"
    ,"      Family and consumer sciences and human sciences"
    ,"      Geosciences, atmospheric sciences, and ocean sciences"
    ,"      Mathematics and statistics"
    ,"      Multidisciplinary and interdisciplinary studies"
    ,"      Natural resources and conservation"
    ,"      Neurobiology and neuroscience"
    ,"      Physical sciences"
    ,"      Psychology"
    ,"      Social sciences"
    ,"    Engineering"
    ,"      Aerospace, aeronautical, and astronautical engineering"
    ,"      Biological, biomedical, and biosystems engineering"
    ,"      Chemical, petroleum, and chemical-related engineering"
    ,"      Civil, environmental, transportation and related engineering fields"
    ,"      Electrical, electronics, communications and computer engineering"
    ,"      Industrial, manufacturing, systems engineering and operations research"
    ,"      Mechanical engineering"
    ,"      Metallurgical, mining, materials and related engineering fields"
    ,"      Other engineering"
    ,"  Health"
    ,"    Clinical medicine"
    ,"    Other health"
)

$ptData = @(
    ,@(2022, 2021, 2020, 2019, 2018, 2017, 2016, 2015, 2014, 2013)
    ,@("1,703", 676, "1,037", 470, 510, "1,106", 579, 594, 662, 747)
    ,@(980, 404, 624, 292, 305, 627, 323, 311, 352, 429)
    ,@(723, 272, 413, 178, 205, 479, 256, 283, 310, 318)
    ,@("1,328", 632, 832, 419, 462, 857, 534, 547, 607, 676)
    ,@(330, 202, 224, 101, 108, 209, 138, 137, 146, 186)
    ,@(772, 405, 582, 302, 334, 623, 375, 388, 436, 455)
    ,@(57, 29, 25, 16, 17, 23, 26, 25, 32, 34)
    ,@(75, 38, 45, 22, 14, 30, 15, 14, 26, 26)
    ,@(37, 15, 15, 12, 6, 19, 13, 13, 13, 17)
    ,@(0, 0, 1, 0, 0, 3, 0, 0, 1, 2)
    ,@(558, 313, 464, 235, 279, 524, 300, 322, 358, 366)
    ,@(45, 10, 32, 17, 18, 24, 21, 14, 6, 10)
    ,@(226, 25, 26, 16, 20, 25, 21, 22, 25, 35)
    ,@(375, 44, 205, 51, 48, 249, 45, 47, 55, 71)
    ,@("1,508", 593, 942, 449, 442, 947, 509, 517, 605, 686)
    ,@(896, 348, 587, 271, 254, 610, 318, 326, 375, 436)
    ,@(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    ,@(82, 22, 58, 12, 11, 47, 24, 21, 36, 49)
    ,@(0, 0, 0, 0, 0, 0, 12, 11, 15, 10)
    ,@(191, 100, 108, 73, 70, 109, 40, 45, 33, 31)
    ,@(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    ,@(51, 6, 29, 1, 7, 39, 12, 4, 13, 9)
    ,@(44, 27, 29, 23, 22, 31, 17, 15, 25, 19)
    ,@(0, 0, 0, 0, 0, 0, 17, 16, 17, 18)
    ,@(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    ,@(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    ,@(148, 37, 97, 35, 33, 102, 28, 16, 17, 33)
    ,@(67, 19, 43, 28, 21, 46, 5, 13, 11, 9)
    ,@(313, 137, 223, 99, 90, 236, 163, 185, 208, 258)
    ,@(612, 245, 355, 178, 188, 337, 191, 191, 230, 250)
    ,@(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    ,@(39, 10, 20, 1, 6, 8, 4, 0, 6, 3)
    ,@(38, 6, 10, 4, 7, 8, 6, 2, 3, 10)
    ,@(141, 49, 78, 38, 48, 59, 47, 41, 55, 54)
    ,@(216, 97, 139, 32, 37, 55, 71, 69, 81, 91)
    ,@(5, 2, 1, 1, 0, 2, 1, 1, 2, 4)
    ,@(100, 57, 67, 35, 23, 26, 43, 52, 49, 43)
    ,@(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    ,@(73, 24, 40, 67, 67, 179, 19, 26, 34, 45)
    ,@(195, 83, 95, 21, 68, 159, 70, 77, 57, 61)
    ,@(140, 57, 58, 17, 22, 84, 21, 23, 15, 16)
    ,@(55, 26, 37, 4, 46, 75, 49, 54, 42, 45)
)


# --- Write row labels (column A) and the year/value grid (columns B:K) ---
for ($r = 0; $r -lt $ptData.Count; $r++) {
    $rowNum = $r + 1
    $wsPartTime.Cells.Item($rowNum, 1).Value = $ptLabels[$r]
    $rowVals = $ptData[$r]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $wsPartTime.Cells.Item($rowNum, $c + 2).Value = $rowVals[$c]
    }
}

# --- Restore the Full-time sheet's selection (no longer the active tab) ---
$wsFullTime.Range("O21").Select()

# --- Make the new Part-time sheet the active tab, with the full data range selected ---
$wsPartTime.Activate()
$wsPartTime.Range("A1:K43").Select()
